$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "A: ['Discount'] > 0"
$ws.Range("D2").Value = 27.1098
$ws.Range("E2").Value = 0.7175
$ws.Range("F2").Value = 20.759

$ws.Range("C3").Value = "B: ['Discount'] == 0"
$ws.Range("D3").Value = 21.7153
$ws.Range("E3").Value = 0.4824
$ws.Range("F3").Value = 17.5008
